$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Column A: timestamp stored as a plain (inline) string, not a date
$ws.Cells.Item($row, 1).Value = "2025-04-29 11:55:13"

# Column B: numeric metric value
$ws.Cells.Item($row, 2).Value = 270
